# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$CellRef, [string]$NewValue)
    $cell = $Sheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "27.449.55"
$ws.Range("E2").Value = "  +6.47%  "
$ws.Range("D3").Value = "1.811.60"
$ws.Range("E3").Value = "  +6.23%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws "D5" "345.19"
$ws.Range("E5").Value = "  +4.24%  "
Set-TextValue $ws "D6" "0.9995"
$ws.Range("E6").Value = "  -0.01%  "
Set-TextValue $ws "D7" "0.3839"
$ws.Range("E7").Value = "  +4.17%  "
Set-TextValue $ws "D8" "50.13"
$ws.Range("E8").Value = "  +4.10%  "
Set-TextValue $ws "D9" "0.3517"
$ws.Range("E9").Value = "  +6.25%  "
Set-TextValue $ws "D10" "1.236"
$ws.Range("E10").Value = "  +5.78%  "
Set-TextValue $ws "D11" "0.07767"
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("E12").Value = "  +0.07%  "
Set-TextValue $ws "D13" "22.56"
$ws.Range("E13").Value = "  +12.80%  "
Set-TextValue $ws "D14" "6.629"
$ws.Range("E14").Value = "  +7.07%  "
Set-TextValue $ws "D15" "7.224"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("D16").Value = "1.810.82"
$ws.Range("E16").Value = "  +6.32%  "
Set-TextValue $ws "D17" "0.00001126"
$ws.Range("E17").Value = "  +5.56%  "
Set-TextValue $ws "D18" "0.06760"
$ws.Range("E18").Value = "  +2.23%  "
Set-TextValue $ws "D19" "86.93"
$ws.Range("E19").Value = "  +7.20%  "
Set-TextValue $ws "D20" "0.9996"
$ws.Range("E20").Value = "  +0.05%  "
Set-TextValue $ws "D21" "17.81"
$ws.Range("E21").Value = "  +10.16%  "
$ws.Range("E22").Value = "  +8.02%  "
Set-TextValue $ws "D23" "13.18"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "27.432.29"
$ws.Range("E24").Value = "  +6.54%  "
Set-TextValue $ws "D25" "2.468"
$ws.Range("E25").Value = "  +0.08%  "
Set-TextValue $ws "D26" "2.680"
$ws.Range("E26").Value = "  +8.16%  "
$ws.Range("E27").Value = "  +16.13%  "
Set-TextValue $ws "D28" "1.499"
$ws.Range("E28").Value = "  +15.50%  "
Set-TextValue $ws "D29" "154.35"
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("D30").Value = "2.012.96"
$ws.Range("E30").Value = "  +6.62%  "
Set-TextValue $ws "D31" "136.86"
$ws.Range("E31").Value = "  +7.26%  "
Set-TextValue $ws "D32" "6.387"
$ws.Range("E32").Value = "  +7.46%  "
Set-TextValue $ws "D33" "4.078"
$ws.Range("E33").Value = "  -0.48%  "
Set-TextValue $ws "D34" "13.91"
$ws.Range("E34").Value = "  +8.38%  "
Set-TextValue $ws "D35" "0.08818"
$ws.Range("E35").Value = "  +3.72%  "
Set-TextValue $ws "D36" "1.722"
$ws.Range("E36").Value = "  +1.68%  "
Set-TextValue $ws "D37" "5.648"
$ws.Range("E37").Value = "  +5.88%  "
Set-TextValue $ws "D38" "0.7104"
$ws.Range("E38").Value = "  +16.26%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D39" "0.02425"
$ws.Range("E39").Value = "  +7.90%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D40" "0.2270"
$ws.Range("E40").Value = "  +7.17%  "
Set-TextValue $ws "D41" "0.06534"
$ws.Range("E41").Value = "  +5.38%  "
Set-TextValue $ws "D42" "9.017"
$ws.Range("E42").Value = "  +5.77%  "
Set-TextValue $ws "D43" "1.290"
$ws.Range("E43").Value = "  +1.26%  "
Set-TextValue $ws "D44" "14.99"
$ws.Range("E44").Value = "  +0.67%  "
Set-TextValue $ws "D45" "0.6613"
$ws.Range("E45").Value = "  +13.44%  "
Set-TextValue $ws "D46" "0.9995"
$ws.Range("E46").Value = "  +0.02%  "
Set-TextValue $ws "D47" "3.978"
$ws.Range("E47").Value = "  +3.51%  "
Set-TextValue $ws "D48" "2.189"
$ws.Range("E48").Value = "  +9.38%  "
Set-TextValue $ws "D49" "133.14"
$ws.Range("E49").Value = "  +5.20%  "
Set-TextValue $ws "D50" "0.07363"
$ws.Range("E50").Value = "  +2.07%  "
Set-TextValue $ws "D51" "80.73"
$ws.Range("E51").Value = "  +5.66%  "
